$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '92.265.07'
$ws.Range("E2").Value = '  +1.78%  '
$ws.Range("D3").Value = '3.111.61'
$ws.Range("E3").Value = '  -2.59%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.65'
$ws.Range("E5").Value = '  -0.85%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '616.01'
$ws.Range("E6").Value = '  -0.46%  '
$ws.Range("E7").Value = '  -1.88%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.391'
$ws.Range("E8").Value = '  +5.06%  '
$ws.Range("E9").Value = '  -0.06%  '
$ws.Range("D10").Value = '3.106.93'
$ws.Range("E10").Value = '  -2.38%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.738'
$ws.Range("E11").Value = '  -0.13%  '
$ws.Range("E12").Value = '  -1.13%  '
$ws.Range("E13").Value = '  +0.43%  '
$ws.Range("D14").Value = '92.278.99'
$ws.Range("E14").Value = '  +1.65%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '34.37'
$ws.Range("E15").Value = '  -2.59%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.45'
$ws.Range("E16").Value = '  -2.27%  '
$ws.Range("D17").Value = '3.693.98'
$ws.Range("E17").Value = '  -1.53%  '
$ws.Range("D18").Value = '3.112.38'
$ws.Range("E18").Value = '  -1.27%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.73'
$ws.Range("E19").Value = '  +0.39%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.67'
$ws.Range("E20").Value = '  -3.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.79'
$ws.Range("E21").Value = '  -3.87%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.39'
$ws.Range("E22").Value = '  +2.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '446.92'
$ws.Range("E23").Value = '  -1.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000197'
$ws.Range("E24").Value = '  -2.70%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.78'
$ws.Range("E25").Value = '  +0.46%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '86.77'
$ws.Range("E26").Value = '  -2.39%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.85'
$ws.Range("E27").Value = '  -1.12%  '
$ws.Range("D28").Value = '3.276.26'
$ws.Range("E28").Value = '  -1.64%  '
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("E30").Value = '  -5.28%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.232'
$ws.Range("E31").Value = '  -1.39%  '
$ws.Range("E32").Value = '  -1.11%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '9.14'
$ws.Range("E33").Value = '  -2.55%  '
$ws.Range("E34").Value = '  -0.60%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.90'
$ws.Range("E35").Value = '  +2.66%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.159'
$ws.Range("E36").Value = '  -7.65%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '26.16'
$ws.Range("E37").Value = '  -2.97%  '
$ws.Range("E38").Value = '  -3.81%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.87'
$ws.Range("E39").Value = '  +0.96%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '484.99'
$ws.Range("E40").Value = '  -5.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.30'
$ws.Range("E41").Value = '  -3.80%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '23.88'
$ws.Range("E42").Value = '  +8.22%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.434'
$ws.Range("E43").Value = '  -4.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.31'
$ws.Range("E44").Value = '  -4.06%  '
$ws.Range("E45").Value = '  +0.06%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '162.47'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.90'
$ws.Range("E47").Value = '  -2.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.692'
$ws.Range("E48").Value = '  -5.55%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.40'
$ws.Range("E49").Value = '  +0.71%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0337'
$ws.Range("E50").Value = '  +4.07%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.42'
$ws.Range("E51").Value = '  -1.19%  '
